$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Appliances"
$ws.Range("B3").Value = "Space cooling"
$ws.Range("B4").Value = "Space heating"
$ws.Range("B5").Value = "Domestic hot water"
$ws.Range("B6").Value = "Ventilation"
